# Updated cryptos list values (prices + 1h volume change) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.996.31'
$ws.Range('E2').Value = '  -3.23%  '

$ws.Range('D3').Value = '1.796.34'
$ws.Range('E3').Value = '  -3.47%  '

$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('D5').Value = '''307.45'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.97%  '

$ws.Range('E6').Value = '  +0.16%  '

$ws.Range('D7').Value = '''0.4187'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -3.47%  '

$ws.Range('D8').Value = '''0.3568'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -4.05%  '

$ws.Range('D9').Value = '''0.07088'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -4.10%  '

$ws.Range('D10').Value = '''0.8438'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -4.54%  '

$ws.Range('D11').Value = '''20.13'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -5.19%  '

$ws.Range('D12').Value = '1.765.36'
$ws.Range('E12').Value = '  -11.56%  '

$ws.Range('D13').Value = '''5.283'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.76%  '

$ws.Range('D14').Value = '''6.349'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -4.36%  '

$ws.Range('D15').Value = '''0.06754'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.95%  '

$ws.Range('D16').Value = '''1.002'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.09%  '

$ws.Range('D17').Value = '''79.64'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.67%  '

$ws.Range('D18').Value = '''0.000008683'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -4.65%  '

$ws.Range('D20').Value = '''14.98'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -4.28%  '

$ws.Range('D21').Value = '26.852.03'
$ws.Range('E21').Value = '  -5.30%  '

$ws.Range('D22').Value = '''5.049'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.48%  '

$ws.Range('D23').Value = '''10.93'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.66%  '

$ws.Range('D24').Value = '1.990.73'
$ws.Range('E24').Value = '  -6.43%  '

$ws.Range('D25').Value = '''1.933'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.57%  '

$ws.Range('E26').Value = '  -1.63%  '

$ws.Range('E27').Value = '  -5.61%  '

$ws.Range('D28').Value = '''4.995'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -6.91%  '

$ws.Range('D29').Value = '''112.92'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.17%  '

$ws.Range('E30').Value = '  -12.54%  '

$ws.Range('D31').Value = '''0.08966'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.11%  '

$ws.Range('D32').Value = '''0.7178'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -9.53%  '

$ws.Range('E33').Value = '  -4.29%  '

$ws.Range('D34').Value = '''4.287'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -7.84%  '

$ws.Range('D36').Value = '''1.075'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -9.10%  '

$ws.Range('D37').Value = '''1.072'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.74%  '

$ws.Range('D38').Value = '''0.01899'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.77%  '

$ws.Range('D39').Value = '''0.05099'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -6.85%  '

$ws.Range('E40').Value = '  -4.41%  '

$ws.Range('D41').Value = '''0.4933'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -5.58%  '

$ws.Range('E42').Value = '  -9.72%  '

$ws.Range('D43').Value = '''5.946'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -12.85%  '

$ws.Range('E44').Value = '  -8.35%  '

$ws.Range('D45').Value = '''104.36'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.72%  '

$ws.Range('D46').Value = '''10.20'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.90%  '

$ws.Range('E47').Value = '  +0.22%  '

$ws.Range('E48').Value = '  -4.29%  '

$ws.Range('D49').Value = '''0.4505'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -6.55%  '

$ws.Range('D50').Value = '''1.597'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -5.09%  '

$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '''1.684'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -9.33%  '

